$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.101.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.917.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.47%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "599.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.92%  "

$ws.Range("E9").Value = "  +1.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.915.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.428"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.96%  "

$ws.Range("E12").Value = "  -0.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.453.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.004.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000191"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.912.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.26%  "

$ws.Range("E19").Value = "  -2.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.065.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000109"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.14%  "

$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "506.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.72%  "

$ws.Range("E34").Value = "  -0.19%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.114"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.09%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "180.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.61%  "

$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0943"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.69%  "

$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("E48").Value = "  -1.17%  "

$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.663"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
